$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "If there are N cores running the parallelisable part, this means that the fastest the parallelisable part can be run at is",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If there are N cores running the parallelisable part, this means that the fastest the parallelisable part can be run at is:",
    2)
